$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "58.264.23"
$ws.Cells.Item(2, 5).Value = "  +0.02%  "
$ws.Cells.Item(3, 4).Value = "2.522.34"
$ws.Cells.Item(3, 5).Value = "  +1.99%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.20%  "
$ws.Cells.Item(5, 4).Value = "'521.50"
$ws.Cells.Item(5, 5).Value = "  +0.11%  "
$ws.Cells.Item(6, 4).Value = "'133.03"
$ws.Cells.Item(6, 5).Value = "  -1.00%  "
$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 5).Value = "  -0.04%  "
$ws.Cells.Item(8, 5).Value = "  +0.52%  "
$ws.Cells.Item(9, 4).Value = "2.520.53"
$ws.Cells.Item(9, 5).Value = "  +1.51%  "
$ws.Cells.Item(10, 4).Value = "'0.0978"
$ws.Cells.Item(10, 5).Value = "  -0.50%  "
$ws.Cells.Item(11, 5).Value = "  -1.45%  "
$ws.Cells.Item(12, 4).Value = "'5.15"
$ws.Cells.Item(12, 5).Value = "  -3.31%  "
$ws.Cells.Item(14, 4).Value = "2.969.53"
$ws.Cells.Item(14, 5).Value = "  +1.99%  "
$ws.Cells.Item(15, 4).Value = "58.313.24"
$ws.Cells.Item(15, 5).Value = "  +0.27%  "
$ws.Cells.Item(16, 4).Value = "'22.10"
$ws.Cells.Item(16, 5).Value = "  -0.36%  "
$ws.Cells.Item(17, 5).Value = "  -0.35%  "
$ws.Cells.Item(18, 4).Value = "2.512.08"
$ws.Cells.Item(18, 5).Value = "  +1.31%  "
$ws.Cells.Item(19, 4).Value = "'10.66"
$ws.Cells.Item(19, 5).Value = "  -0.12%  "
$ws.Cells.Item(20, 4).Value = "'321.65"
$ws.Cells.Item(20, 5).Value = "  +0.39%  "
$ws.Cells.Item(21, 5).Value = "  -0.85%  "
$ws.Cells.Item(22, 4).Value = "'6.16"
$ws.Cells.Item(22, 5).Value = "  +7.51%  "
$ws.Cells.Item(23, 4).Value = "'0.999"
$ws.Cells.Item(23, 5).Value = "  -0.01%  "
$ws.Cells.Item(24, 4).Value = "'64.49"
$ws.Cells.Item(24, 5).Value = "  +0.02%  "
$ws.Cells.Item(25, 4).Value = "'0.407"
$ws.Cells.Item(25, 5).Value = "  -0.83%  "
$ws.Cells.Item(26, 4).Value = "'0.998"
$ws.Cells.Item(26, 5).Value = "  +0.16%  "
$ws.Cells.Item(27, 5).Value = "  -0.45%  "
$ws.Cells.Item(28, 4).Value = "'7.38"
$ws.Cells.Item(28, 5).Value = "  -0.03%  "
$ws.Cells.Item(29, 4).Value = "0.0₃0753"
$ws.Cells.Item(29, 5).Value = "  +0.42%  "
$ws.Cells.Item(30, 5).Value = "  +1.50%  "
$ws.Cells.Item(31, 4).Value = "'167.43"
$ws.Cells.Item(31, 5).Value = "  -1.21%  "
$ws.Cells.Item(32, 5).Value = "  +0.84%  "
$ws.Cells.Item(33, 4).Value = "'6.31"
$ws.Cells.Item(33, 5).Value = "  +0.25%  "
$ws.Cells.Item(34, 4).Value = "'0.997"
$ws.Cells.Item(34, 5).Value = "  -0.09%  "
$ws.Cells.Item(35, 4).Value = "'0.997"
$ws.Cells.Item(35, 5).Value = "  +0.10%  "
$ws.Cells.Item(36, 4).Value = "'18.18"
$ws.Cells.Item(36, 5).Value = "  +0.45%  "
$ws.Cells.Item(37, 5).Value = "  -5.25%  "
$ws.Cells.Item(38, 4).Value = "'3.92"
$ws.Cells.Item(38, 5).Value = "  -2.21%  "
$ws.Cells.Item(39, 5).Value = "  +0.72%  "
$ws.Cells.Item(40, 4).Value = "'36.42"
$ws.Cells.Item(40, 5).Value = "  -0.44%  "
$ws.Cells.Item(41, 4).Value = "'0.771"
$ws.Cells.Item(41, 5).Value = "  -3.69%  "
$ws.Cells.Item(42, 2).Value = "Bittensor"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(42, 4).Value = "'276.59"
$ws.Cells.Item(42, 5).Value = "  +1.02%  "
$ws.Cells.Item(43, 2).Value = "Filecoin"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(43, 4).Value = "'3.44"
$ws.Cells.Item(43, 5).Value = "  -0.25%  "
$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(44, 4).Value = "'4.99"
$ws.Cells.Item(44, 5).Value = "  -3.36%  "
$ws.Cells.Item(45, 2).Value = "Aave"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(45, 4).Value = "'129.24"
$ws.Cells.Item(45, 5).Value = "  +4.05%  "
$ws.Cells.Item(46, 4).Value = "'0.598"
$ws.Cells.Item(46, 5).Value = "  +0.33%  "
$ws.Cells.Item(47, 5).Value = "  +1.37%  "
$ws.Cells.Item(48, 4).Value = "'0.0501"
$ws.Cells.Item(48, 5).Value = "  +2.27%  "
$ws.Cells.Item(49, 4).Value = "'17.71"
$ws.Cells.Item(49, 5).Value = "  -0.46%  "
$ws.Cells.Item(50, 5).Value = "  +0.13%  "
$ws.Cells.Item(51, 4).Value = "'16.90"
$ws.Cells.Item(51, 5).Value = "  -0.87%  "
